$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values for row 4 (Persons Receiving Social Package)
$ws.Range("E4").Value = 5189
$ws.Range("F4").Value = 5283
$ws.Range("G4").Value = 5280
$ws.Range("H4").Value = 5231

# Update the selection/active cell on the sheet (A3 becomes active cell, sqref A3)
$ws.Range("A3").Select()
